$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = 42613.75717592592
$ws.Range("B11").Value = 34
$ws.Range("C11").Value = 59
$ws.Range("D11").Value = 39
$ws.Range("E11").Value = 59
$ws.Range("F11").Value = 14
$ws.Range("G11").Value = 17593
$ws.Range("H11").Value = 16615
$ws.Range("I11").Value = 2888
$ws.Range("J11").Value = 371
$ws.Range("K11").Value = 247
$ws.Range("L11").Value = 41
$ws.Range("M11").Value = 7
$ws.Range("N11").Value = "Noun"

$ws.Range("A12").Value = 42613.88434027778
$ws.Range("B12").Value = 20
$ws.Range("C12").Value = 60
$ws.Range("D12").Value = 38
$ws.Range("E12").Value = 60
$ws.Range("F12").Value = 36
$ws.Range("G12").Value = 11283
$ws.Range("H12").Value = 9814
$ws.Range("I12").Value = 1705
$ws.Range("J12").Value = 216
$ws.Range("K12").Value = 136
$ws.Range("L12").Value = 23
$ws.Range("M12").Value = 13
$ws.Range("N12").Value = "Noun"

$ws.Range("A13").Value = 42614.883101851854
$ws.Range("B13").Value = 8
$ws.Range("C13").Value = 50
$ws.Range("D13").Value = 45
$ws.Range("E13").Value = 50
$ws.Range("F13").Value = 25
$ws.Range("G13").Value = 21226
$ws.Range("H13").Value = 21847
$ws.Range("I13").Value = 3720
$ws.Range("J13").Value = 424
$ws.Range("K13").Value = 379
$ws.Range("L13").Value = 43
$ws.Range("M13").Value = 15
$ws.Range("N13").Value = "Noun"

$ws.Range("A14").Value = 42615.8840162037
$ws.Range("B14").Value = 16
$ws.Range("C14").Value = 51
$ws.Range("D14").Value = 42
$ws.Range("E14").Value = 51
$ws.Range("F14").Value = 25
$ws.Range("G14").Value = 16845
$ws.Range("H14").Value = 19094
$ws.Range("I14").Value = 3224
$ws.Range("J14").Value = 371
$ws.Range("K14").Value = 302
$ws.Range("L14").Value = 47
$ws.Range("M14").Value = 16
$ws.Range("N14").Value = "Noun"
